$p = $ppt.ActivePresentation

# --- Slide 1: title subtitle box - normalize the "IIR.C11" / ". Probabilistic "
#     run split into one run so it reads "IIR.C11. Probabilistic " ---
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$fullText1 = $tr1.Text
$markerLen = "IIR.C11. Probabilistic ".Length
$startPos = $fullText1.IndexOf("IIR.C11") + 1
$sub1 = $tr1.Characters($startPos, $markerLen)
$sub1.Text = "IIR.C11. Probabilistic "

# --- Slide 36: update exercise title "Bài tập" -> "Bài tập 5.1" ---
$s36 = $p.Slides.Item(36)
$shp36 = $s36.Shapes.Item(1)
$tr36 = $shp36.TextFrame.TextRange
$fullText36 = $tr36.Text
$tapPos = $fullText36.IndexOf("tập") + 1
$sub36 = $tr36.Characters($tapPos, "tập".Length)
$sub36.Text = "tập 5.1"
